$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.431.69'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.03%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.566.19'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.39%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("E5").Value = '  -0.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '285.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.43%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3634'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.27%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '48.16'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.32%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3320'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.44%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.121'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.28%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07394'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.07%  '
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.72'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.49%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.943'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.893'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.99%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.568.52'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.68%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001102'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.89%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '87.46'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06717'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.67%  '
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.381'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.25%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.21'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.56%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.420.21'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.376'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.551'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '150.89'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.49%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.41'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.998'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.78%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.96'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.53%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.740.67'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.53%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.022'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.29%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.996'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.14%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.086'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.80%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.709'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08234'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.57%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02408'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.85%  '
$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2230'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.14%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06409'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.66%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.293'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.96%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.367'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.58%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6254'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.45%  '
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.19'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.78%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.82'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.90%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6036'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.86%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.742'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.69%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.026'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.35%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '123.37'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.211'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07200'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.65%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '75.76'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.36%  '
